# Auto-generated Excel COM-interop script
# Applies cell-value corrections to the Seraph_Profits workbook per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1035.3334
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H15").Value = 422.70834
$ws.Range("I15").Value = 422.70834
$ws.Range("K15").Value = 1268.12502
$ws.Range("M15").Value = -1099.12502

$ws.Range("H86").Value = 1200.75
$ws.Range("I86").Value = 1201
$ws.Range("J86").Value = 1200
$ws.Range("K86").Value = 1201
$ws.Range("L86").Value = 1200
$ws.Range("M86").Value = -78
$ws.Range("N86").Value = -3446

$ws.Range("H87").Value = 75998.5
$ws.Range("J87").Value = 75998.5
$ws.Range("L87").Value = 75998.5
$ws.Range("N87").Value = -78494.5

$ws.Range("H89").Value = 1200.75
$ws.Range("I89").Value = 1201
$ws.Range("J89").Value = 1200
$ws.Range("K89").Value = 6005
$ws.Range("L89").Value = 6000
$ws.Range("M89").Value = -389
$ws.Range("N89").Value = -17232

$ws.Range("H90").Value = 75998.5
$ws.Range("J90").Value = 75998.5
$ws.Range("L90").Value = 227995.5
$ws.Range("N90").Value = -240475.5

$ws.Range("H113").Value = 2347
$ws.Range("I113").Value = 3247.5
$ws.Range("J113").Value = 1446.5
$ws.Range("K113").Value = 3247.5
$ws.Range("L113").Value = 1446.5
$ws.Range("M113").Value = 6.5
$ws.Range("N113").Value = -7954.5

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H137").Value = 2164.4
$ws.Range("I137").Value = 2107.3333
$ws.Range("K137").Value = 6321.999899999999
$ws.Range("M137").Value = -3771.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 3000
$ws.Range("J30").Value = 3000
$ws.Range("L30").Value = 3000
$ws.Range("N30").Value = -3250

$ws.Range("H40").Value = 49333.332
$ws.Range("J40").Value = 49333.332
$ws.Range("L40").Value = 49333.332
$ws.Range("N40").Value = -49863.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3992.6956
$ws.Range("I31").Value = 2697.2
$ws.Range("J31").Value = 6421.75
$ws.Range("K31").Value = 2697.2
$ws.Range("L31").Value = 6421.75
$ws.Range("M31").Value = -2402.2
$ws.Range("N31").Value = -7011.75

$ws.Range("H34").Value = 3992.6956
$ws.Range("I34").Value = 2697.2
$ws.Range("J34").Value = 6421.75
$ws.Range("K34").Value = 2697.2
$ws.Range("L34").Value = 6421.75
$ws.Range("M34").Value = -2495.2
$ws.Range("N34").Value = -6825.75

$ws.Range("H41").Value = 8902.333000000001
$ws.Range("I41").Value = 4799
$ws.Range("J41").Value = 29419
$ws.Range("K41").Value = 4799
$ws.Range("L41").Value = 29419
$ws.Range("M41").Value = -4371
$ws.Range("N41").Value = -30275

$ws.Range("H50").Value = 18873.875
$ws.Range("I50").Value = 6999.5
$ws.Range("J50").Value = 22832
$ws.Range("K50").Value = 6999.5
$ws.Range("L50").Value = 22832
$ws.Range("M50").Value = -6374.5
$ws.Range("N50").Value = -24082

$ws.Range("H51").Value = 29998
$ws.Range("J51").Value = 29998
$ws.Range("L51").Value = 29998
$ws.Range("N51").Value = -31470

$ws.Range("H59").Value = 48973.25
$ws.Range("J59").Value = 59949
$ws.Range("L59").Value = 59949
$ws.Range("N59").Value = -62239

$ws.Range("H60").Value = 13264.363
$ws.Range("J60").Value = 49988
$ws.Range("L60").Value = 49988
$ws.Range("N60").Value = -51010

$ws.Range("H61").Value = 29998
$ws.Range("J61").Value = 29998
$ws.Range("L61").Value = 29998
$ws.Range("N61").Value = -30694

$ws.Range("H68").Value = 37498.5
$ws.Range("J68").Value = 37498.5
$ws.Range("L68").Value = 37498.5
$ws.Range("N68").Value = -38996.5

$ws.Range("H71").Value = 37498.5
$ws.Range("J71").Value = 37498.5
$ws.Range("L71").Value = 112495.5
$ws.Range("N71").Value = -119983.5

$ws.Range("H74").Value = 48247.75
$ws.Range("J74").Value = 48247.75
$ws.Range("L74").Value = 48247.75
$ws.Range("N74").Value = -49995.75

$ws.Range("H77").Value = 48247.75
$ws.Range("J77").Value = 48247.75
$ws.Range("L77").Value = 144743.25
$ws.Range("N77").Value = -153479.25

$ws.Range("H122").Value = 980.375
$ws.Range("I122").Value = 909
$ws.Range("J122").Value = 1099.3334
$ws.Range("K122").Value = 2727
$ws.Range("L122").Value = 3298.0002
$ws.Range("M122").Value = -277
$ws.Range("N122").Value = -8198.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 80.2
$ws.Range("I6").Value = 80.2
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 240.6
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -127.6

$ws.Range("H81").Value = 2995
$ws.Range("J81").Value = 3490
$ws.Range("L81").Value = 10470
$ws.Range("N81").Value = -12716

$ws.Range("H84").Value = 2995
$ws.Range("J84").Value = 3490
$ws.Range("L84").Value = 31410
$ws.Range("N84").Value = -42642

$ws.Range("H138").Value = 4189.75
$ws.Range("I138").Value = 3930
$ws.Range("J138").Value = 4969
$ws.Range("K138").Value = 11790
$ws.Range("L138").Value = 14907
$ws.Range("M138").Value = -6650
$ws.Range("N138").Value = -25187

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3800.625
$ws.Range("I80").Value = 2279.8
$ws.Range("J80").Value = 6335.3335
$ws.Range("K80").Value = 2279.8
$ws.Range("L80").Value = 6335.3335
$ws.Range("M80").Value = -1281.8
$ws.Range("N80").Value = -8331.333500000001

$ws.Range("H83").Value = 3800.625
$ws.Range("I83").Value = 2279.8
$ws.Range("J83").Value = 6335.3335
$ws.Range("K83").Value = 11399
$ws.Range("L83").Value = 31676.6675
$ws.Range("M83").Value = -6407
$ws.Range("N83").Value = -41660.6675

$ws.Range("H113").Value = 44712.855
$ws.Range("I113").Value = 35666
$ws.Range("J113").Value = 51498
$ws.Range("K113").Value = 35666
$ws.Range("L113").Value = 51498
$ws.Range("M113").Value = -33496
$ws.Range("N113").Value = -55838

$ws.Range("H132").Value = 2033.1666
$ws.Range("I132").Value = 1639.8
$ws.Range("K132").Value = 4919.4
$ws.Range("M132").Value = -2389.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 20000
$ws.Range("J26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("N26").Value = -20590

$ws.Range("H46").Value = 204599.6
$ws.Range("I46").Value = 6999.5
$ws.Range("J46").Value = 336333
$ws.Range("K46").Value = 6999.5
$ws.Range("L46").Value = 336333
$ws.Range("M46").Value = -6811.5
$ws.Range("N46").Value = -336709

$ws.Range("H82").Value = 859.8570999999999
$ws.Range("I82").Value = 694.1
$ws.Range("J82").Value = 1274.25
$ws.Range("K82").Value = 694.1
$ws.Range("L82").Value = 1274.25
$ws.Range("M82").Value = -333.1
$ws.Range("N82").Value = -1996.25

$ws.Range("H85").Value = 859.8570999999999
$ws.Range("I85").Value = 694.1
$ws.Range("J85").Value = 1274.25
$ws.Range("K85").Value = 694.1
$ws.Range("L85").Value = 1274.25
$ws.Range("M85").Value = 553.9
$ws.Range("N85").Value = -3770.25

$ws.Range("H136").Value = 6798.6665
$ws.Range("I136").Value = 6608.727
$ws.Range("K136").Value = 19826.181
$ws.Range("M136").Value = -17276.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 47665.668
$ws.Range("J64").Value = 47665.668
$ws.Range("L64").Value = 47665.668
$ws.Range("N64").Value = -48161.668

$ws.Range("H67").Value = 47665.668
$ws.Range("J67").Value = 47665.668
$ws.Range("L67").Value = 47665.668
$ws.Range("N67").Value = -49381.668

